$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("D3").Value = 10.43
$ws.Range("E3").Value = 10.88
$ws.Range("F3").Value = 10.38
$ws.Range("J3").Value = 9.5

$ws.Range("C4").Value = 9.57
$ws.Range("E4").Value = 10.59
$ws.Range("F4").Value = 10.09

$ws.Range("C5").Value = 9.119999999999999
$ws.Range("D5").Value = 9.41
$ws.Range("F5").Value = 10.29
$ws.Range("H5").Value = 8.880000000000001

$ws.Range("C6").Value = 9.619999999999999
$ws.Range("D6").Value = 9.91
$ws.Range("E6").Value = 9.710000000000001
$ws.Range("G6").Value = 10.43

$ws.Range("F7").Value = 9.57
$ws.Range("H7").Value = 9.69

$ws.Range("E8").Value = 11.12
$ws.Range("G8").Value = 10.31
$ws.Range("I8").Value = 7.33

$ws.Range("H9").Value = 12.67

$ws.Range("C10").Value = 10.5
